# "Save vacations history with deleted user"
# Row 2 is rewritten with a new (surviving) vacation record, and the old
# rows 3-5 (which belonged to users no longer present) are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite row 2 ---------------------------------------------------
$ws.Range("A2").Value = "test3"

# Columns B, D and E hold date-looking text (e.g. "2023-07-10"); format
# them as Text first so Excel keeps the literal string instead of
# auto-converting it to a date serial number.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2023-07-10"

$ws.Range("C2").Value = "wypoczynkowy"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2023-07-14"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2023-07-16"

$ws.Range("F2").Value = 2
$ws.Range("G2").Value = "Tak"

# --- Remove the old rows 3-5 (deleted user's vacation history) -------
$ws.Range("A3:G5").EntireRow.Delete()
